$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.655.37"
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("D3").Value = "2.006.36"
$ws.Range("E3").Value = "  -4.87%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  +0.66%  "
$ws.Range("D5").Value = "'331.84"
$ws.Range("E5").Value = "  -3.88%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "'0.5026"
$ws.Range("E7").Value = "  -3.95%  "
$ws.Range("D8").Value = "'0.4262"
$ws.Range("E8").Value = "  -4.11%  "
$ws.Range("D9").Value = "'54.74"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'0.09164"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").Value = "'1.127"
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("D12").Value = "'23.55"
$ws.Range("D13").Value = "'8.141"
$ws.Range("E13").Value = "  -6.40%  "
$ws.Range("D14").Value = "2.010.46"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "'6.550"
$ws.Range("E15").Value = "  -5.73%  "
$ws.Range("D16").Value = "'95.41"
$ws.Range("E16").Value = "  -6.42%  "
$ws.Range("D17").Value = "'1.014"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'0.00001124"
$ws.Range("E18").Value = "  -3.43%  "
$ws.Range("D19").Value = "'0.06676"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "'19.89"
$ws.Range("E20").Value = "  -6.24%  "
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "'5.990"
$ws.Range("E22").Value = "  -5.49%  "
$ws.Range("D23").Value = "29.654.05"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").Value = "'12.07"
$ws.Range("E24").Value = "  -4.65%  "
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'159.38"
$ws.Range("E26").Value = "  -2.22%  "
$ws.Range("D27").Value = "'20.83"
$ws.Range("E27").Value = "  -5.54%  "
$ws.Range("D28").Value = "'6.422"
$ws.Range("E28").Value = "  -5.75%  "
$ws.Range("D29").Value = "'2.332"
$ws.Range("E29").Value = "  -7.98%  "
$ws.Range("D30").Value = "'128.90"
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("D31").Value = "'1.066"
$ws.Range("E31").Value = "  -7.68%  "
$ws.Range("E32").Value = "  -8.92%  "
$ws.Range("D33").Value = "'0.09963"
$ws.Range("E33").Value = "  -5.56%  "
$ws.Range("D34").Value = "'5.863"
$ws.Range("E34").Value = "  -6.48%  "
$ws.Range("D35").Value = "'3.809"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").Value = "'9.567"
$ws.Range("E36").Value = "  -8.00%  "
$ws.Range("D37").Value = "'0.02480"
$ws.Range("E37").Value = "  -5.48%  "
$ws.Range("D38").Value = "'1.321"
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("D39").Value = "'0.06399"
$ws.Range("E39").Value = "  -5.81%  "
$ws.Range("D40").Value = "'0.6603"
$ws.Range("E40").Value = "  -6.41%  "
$ws.Range("E41").Value = "  -6.27%  "
$ws.Range("D42").Value = "'0.2075"
$ws.Range("E42").Value = "  -6.81%  "
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("E44").Value = "  -6.99%  "
$ws.Range("D45").Value = "'13.62"
$ws.Range("E45").Value = "  -6.30%  "
$ws.Range("D46").Value = "'2.219"
$ws.Range("E46").Value = "  -6.02%  "
$ws.Range("D47").Value = "'1.295"
$ws.Range("E47").Value = "  -5.39%  "
$ws.Range("D48").Value = "'3.536"
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("D49").Value = "'0.07005"
$ws.Range("E49").Value = "  -3.41%  "
$ws.Range("D50").Value = "'0.00000000324"
$ws.Range("E50").Value = "  -6.76%  "
$ws.Range("D51").Value = "'1.133"
$ws.Range("E51").Value = "  -6.01%  "
